# Powerpoint writer: consolidate text run nodes.
# Merge each "word + following space" pair of runs into a single run,
# reducing the number of <a:r> nodes generated for split text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1: Title "A slide" -> merge "A" + " " into a single run "A "
$titleShape = $s.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Characters(1, 2).Text = "A "

# Shape 4: TextBox "Just an image on this side"
# -> merge each word with its trailing space into one run each,
#    leaving the final word ("side") as its own run.
$textBoxShape = $s.Shapes.Item(4)
$textBoxRange = $textBoxShape.TextFrame.TextRange

$textBoxRange.Characters(1, 5).Text = "Just "
$textBoxRange.Characters(6, 3).Text = "an "
$textBoxRange.Characters(9, 6).Text = "image "
$textBoxRange.Characters(15, 3).Text = "on "
$textBoxRange.Characters(18, 5).Text = "this "
